$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 2000
$ws.Range("I69").Value = 1000
$ws.Range("J69").Value = 2500
$ws.Range("K69").Value = 3000
$ws.Range("L69").Value = 7500
$ws.Range("M69").Value = -2126
$ws.Range("N69").Value = -9248

$ws.Range("H72").Value = 2000
$ws.Range("I72").Value = 1000
$ws.Range("J72").Value = 2500
$ws.Range("K72").Value = 9000
$ws.Range("L72").Value = 22500
$ws.Range("M72").Value = -4632
$ws.Range("N72").Value = -31236

$ws.Range("H92").Value = 613.0833
$ws.Range("I92").Value = 613.0833
$ws.Range("K92").Value = 613.0833
$ws.Range("M92").Value = 634.9167

$ws.Range("H137").Value = 1323.6129
$ws.Range("I137").Value = 1294.5186
$ws.Range("J137").Value = 1520
$ws.Range("K137").Value = 3883.5558
$ws.Range("L137").Value = 4560
$ws.Range("M137").Value = -1333.5558
$ws.Range("N137").Value = -9660

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2999.8
$ws.Range("I74").Value = 2499.75
$ws.Range("J74").Value = 5000
$ws.Range("K74").Value = 2499.75
$ws.Range("L74").Value = 5000
$ws.Range("M74").Value = -1625.75
$ws.Range("N74").Value = -6748

$ws.Range("H77").Value = 2999.8
$ws.Range("I77").Value = 2499.75
$ws.Range("J77").Value = 5000
$ws.Range("K77").Value = 12498.75
$ws.Range("L77").Value = 25000
$ws.Range("M77").Value = -8130.75
$ws.Range("N77").Value = -33736

$ws.Range("H92").Value = 149441.5
$ws.Range("J92").Value = 149441.5
$ws.Range("L92").Value = 149441.5
$ws.Range("N92").Value = -154433.5

$ws.Range("H110").Value = 1737.7693
$ws.Range("I110").Value = 1654.6666
$ws.Range("J110").Value = 1924.75
$ws.Range("K110").Value = 1654.6666
$ws.Range("L110").Value = 1924.75
$ws.Range("M110").Value = 390.3334
$ws.Range("N110").Value = -6014.75

$ws.Range("H123").Value = 24287.8
$ws.Range("J123").Value = 24287.8
$ws.Range("L123").Value = 24287.8
$ws.Range("N123").Value = -34087.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H137").Value = 49426.668
$ws.Range("J137").Value = 49426.668
$ws.Range("L137").Value = 49426.668
$ws.Range("N137").Value = -59626.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 389.6087
$ws.Range("I22").Value = 227.58824
$ws.Range("J22").Value = 848.6667
$ws.Range("K22").Value = 227.58824
$ws.Range("L22").Value = 848.6667
$ws.Range("M22").Value = 122.41176
$ws.Range("N22").Value = -1548.6667

$ws.Range("H31").Value = 2274.875
$ws.Range("I31").Value = 1149.8182
$ws.Range("J31").Value = 4750
$ws.Range("K31").Value = 1149.8182
$ws.Range("L31").Value = 4750
$ws.Range("M31").Value = -854.8181999999999
$ws.Range("N31").Value = -5340

$ws.Range("H34").Value = 2274.875
$ws.Range("I34").Value = 1149.8182
$ws.Range("J34").Value = 4750
$ws.Range("K34").Value = 1149.8182
$ws.Range("L34").Value = 4750
$ws.Range("M34").Value = -947.8181999999999
$ws.Range("N34").Value = -5154

$ws.Range("H56").Value = 28546.5
$ws.Range("I56").Value = 28546.5
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 28546.5
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -27701.5
$ws.Range("N56").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2655.3333
$ws.Range("I5").Value = 3286.8572
$ws.Range("J5").Value = 445
$ws.Range("K5").Value = 9860.571599999999
$ws.Range("L5").Value = 1335
$ws.Range("M5").Value = -9748.571599999999
$ws.Range("N5").Value = -1559

$ws.Range("H107").Value = 342.94446
$ws.Range("I107").Value = 391
$ws.Range("J107").Value = 275.66666
$ws.Range("K107").Value = 1173
$ws.Range("L107").Value = 826.9999799999999
$ws.Range("M107").Value = 747
$ws.Range("N107").Value = -4666.99998

$ws.Range("H119").Value = 3881
$ws.Range("I119").Value = 3116.125
$ws.Range("K119").Value = 9348.375
$ws.Range("M119").Value = -4510.375

$ws.Range("H135").Value = 2655.3333
$ws.Range("I135").Value = 3286.8572
$ws.Range("J135").Value = 445
$ws.Range("K135").Value = 29581.7148
$ws.Range("L135").Value = 4005
$ws.Range("M135").Value = -27046.7148
$ws.Range("N135").Value = -9075

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3347.647
$ws.Range("I80").Value = 3408.4614
$ws.Range("J80").Value = 3150
$ws.Range("K80").Value = 3408.4614
$ws.Range("L80").Value = 3150
$ws.Range("M80").Value = -2410.4614
$ws.Range("N80").Value = -5146

$ws.Range("H83").Value = 3347.647
$ws.Range("I83").Value = 3408.4614
$ws.Range("J83").Value = 3150
$ws.Range("K83").Value = 17042.307
$ws.Range("L83").Value = 15750
$ws.Range("M83").Value = -12050.307
$ws.Range("N83").Value = -25734

$ws.Range("H109").Value = 10465.375
$ws.Range("J109").Value = 10465.375
$ws.Range("L109").Value = 10465.375
$ws.Range("N109").Value = -12545.375

$ws.Range("H123").Value = 8617
$ws.Range("J123").Value = 8617
$ws.Range("L123").Value = 8617
$ws.Range("N123").Value = -13517

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 37409944
$ws.Range("I136").Value = 55557224
$ws.Range("K136").Value = 166671672
$ws.Range("M136").Value = -166669122

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 1499.5
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()

$ws.Range("H81").Value = 68173.39999999999
$ws.Range("I81").Value = 84891.586
$ws.Range("J81").Value = 1300.6666
$ws.Range("K81").Value = 169783.172
$ws.Range("L81").Value = 2601.3332
$ws.Range("M81").Value = -168722.172
$ws.Range("N81").Value = -4723.3332

$ws.Range("H84").Value = 68173.39999999999
$ws.Range("I84").Value = 84891.586
$ws.Range("J84").Value = 1300.6666
$ws.Range("K84").Value = 848915.86
$ws.Range("L84").Value = 13006.666
$ws.Range("M84").Value = -843611.86
$ws.Range("N84").Value = -23614.666

$ws.Range("H99").Value = 181666.67
$ws.Range("I99").Value = 260000
$ws.Range("J99").Value = 25000
$ws.Range("K99").Value = 260000
$ws.Range("L99").Value = 25000
$ws.Range("M99").Value = -257005
$ws.Range("N99").Value = -30990

$ws.Range("H123").Value = 22405.318
$ws.Range("J123").Value = 22405.318
$ws.Range("L123").Value = 22405.318
$ws.Range("N123").Value = -32205.318

$ws.Range("H136").Value = 2417.0833
$ws.Range("I136").Value = 2715
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 8145
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -5595
$ws.Range("N136").Value = -11100
